# Restore the "Parameter"/"Value" table: shared-string order was reshuffled
# (column A) and the computed results were refreshed (column B).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "param_TimeStep_starting_index"
$ws.Range("B2").Value = 40
$ws.Range("A3").Value = "param_demand1_op_cost_starting_index"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "param_demand1_inv_cost_starting_index"
$ws.Range("B4").Value = 0
$ws.Range("A5").Value = "param_demand2_inv_cost_starting_index"
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = "param_demand2_op_cost_starting_index"
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = "param_Q_net1_demand2_starting_index"
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = "param_net1_sell_thermal_starting_index"
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = "param_net1_buy_electric_starting_index"
$ws.Range("B9").Value = 36.31986708533334
$ws.Range("A10").Value = "param_net1_sell_electric_starting_index"
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = "param_net1_emissions_starting_index"
$ws.Range("B11").Value = 115.3599238027987
$ws.Range("A12").Value = "param_P_net1_bat2_starting_index"
$ws.Range("B12").Value = 0
$ws.Range("A13").Value = "param_P_net1_heat_pump2_starting_index"
$ws.Range("B13").Value = 19.8888888
$ws.Range("A14").Value = "param_P_net1_bat1_starting_index"
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = "param_P_net1_charging_station1_starting_index"
$ws.Range("B15").Value = 9.93404166666668
$ws.Range("A16").Value = "param_P_to_net1_starting_index"
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = "param_P_net1_demand2_starting_index"
$ws.Range("B17").Value = 0
$ws.Range("A18").Value = "param_P_from_net1_starting_index"
$ws.Range("B18").Value = 82.54515246666668
$ws.Range("A19").Value = "param_P_net1_demand1_starting_index"
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = "param_Q_from_net1_starting_index"
$ws.Range("B20").Value = 288.0609934227724
$ws.Range("A21").Value = "param_Q_to_net1_starting_index"
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = "param_net1_inv_cost_starting_index"
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = "param_P_net1_charging_station2_starting_index"
$ws.Range("B23").Value = 32.8333332
$ws.Range("A24").Value = "param_P_net1_heat_pump1_starting_index"
$ws.Range("B24").Value = 19.8888888
$ws.Range("A25").Value = "param_Q_net1_demand1_starting_index"
$ws.Range("B25").Value = 288.0609934227724
$ws.Range("A26").Value = "param_net1_buy_thermal_starting_index"
$ws.Range("B26").Value = 92.17951789528718
$ws.Range("A27").Value = "param_net2_buy_electric_starting_index"
$ws.Range("B27").Value = 340.4589316647393
$ws.Range("A28").Value = "param_P_net2_bat1_starting_index"
$ws.Range("B28").Value = 0
$ws.Range("A29").Value = "param_net2_inv_cost_starting_index"
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = "param_net2_buy_thermal_starting_index"
$ws.Range("B30").Value = 200
$ws.Range("A31").Value = "param_net2_sell_thermal_starting_index"
$ws.Range("B31").Value = 0
$ws.Range("A32").Value = "param_P_net2_bat2_starting_index"
$ws.Range("B32").Value = 0
$ws.Range("A33").Value = "param_P_net2_charging_station1_starting_index"
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = "param_Q_from_net2_starting_index"
$ws.Range("B34").Value = 1000
$ws.Range("A35").Value = "param_P_net2_heat_pump2_starting_index"
$ws.Range("B35").Value = 0
$ws.Range("A36").Value = "param_P_from_net2_starting_index"
$ws.Range("B36").Value = 851.1473291618483
$ws.Range("A37").Value = "param_P_net2_demand2_starting_index"
$ws.Range("B37").Value = 500
$ws.Range("A38").Value = "param_Q_net2_demand2_starting_index"
$ws.Range("B38").Value = 1000
$ws.Range("A39").Value = "param_net2_sell_electric_starting_index"
$ws.Range("B39").Value = 0
$ws.Range("A40").Value = "param_P_to_net2_starting_index"
$ws.Range("B40").Value = 0
$ws.Range("A41").Value = "param_P_net2_charging_station2_starting_index"
$ws.Range("B41").Value = 0
$ws.Range("A42").Value = "param_Q_to_net2_starting_index"
$ws.Range("B42").Value = 0
$ws.Range("A43").Value = "param_net2_emissions_starting_index"
$ws.Range("B43").Value = 716.6425043306351
$ws.Range("A44").Value = "param_P_net2_heat_pump1_starting_index"
$ws.Range("B44").Value = 0
$ws.Range("A45").Value = "param_P_net2_demand1_starting_index"
$ws.Range("B45").Value = 351.1473291618483
$ws.Range("A46").Value = "param_Q_net2_demand1_starting_index"
$ws.Range("B46").Value = 0
$ws.Range("A47").Value = "param_P_pv1_demand2_starting_index"
$ws.Range("B47").Value = 0
$ws.Range("A48").Value = "param_pv1_op_cost_starting_index"
$ws.Range("B48").Value = 1
$ws.Range("A49").Value = "param_pv1_emissions_starting_index"
$ws.Range("B49").Value = 0.0833334
$ws.Range("A50").Value = "param_P_pv1_demand1_starting_index"
$ws.Range("B50").Value = 0
$ws.Range("A51").Value = "param_P_pv1_bat1_starting_index"
$ws.Range("B51").Value = 0
$ws.Range("A52").Value = "param_P_pv1_charging_station2_starting_index"
$ws.Range("B52").Value = 0.1666668
$ws.Range("A53").Value = "param_P_pv1_bat2_starting_index"
$ws.Range("B53").Value = 0
$ws.Range("A54").Value = "param_P_pv1_net2_starting_index"
$ws.Range("B54").Value = 0
$ws.Range("A55").Value = "param_P_pv1_charging_station1_starting_index"
$ws.Range("B55").Value = 0
$ws.Range("A56").Value = "param_P_from_pv1_starting_index"
$ws.Range("B56").Value = 0.1666668
$ws.Range("A57").Value = "param_P_pv1_heat_pump2_starting_index"
$ws.Range("B57").Value = 0
$ws.Range("A58").Value = "param_P_pv1_heat_pump1_starting_index"
$ws.Range("B58").Value = 0
$ws.Range("A59").Value = "param_P_pv1_net1_starting_index"
$ws.Range("B59").Value = 0
$ws.Range("A60").Value = "param_pv1_inv_cost_starting_index"
$ws.Range("B60").Value = 0
$ws.Range("A61").Value = "param_P_pv2_charging_station1_starting_index"
$ws.Range("B61").Value = 0
$ws.Range("A62").Value = "param_P_pv2_net1_starting_index"
$ws.Range("B62").Value = 0
$ws.Range("A63").Value = "param_P_from_pv2_starting_index"
$ws.Range("B63").Value = 1.8
$ws.Range("A64").Value = "param_P_pv2_bat1_starting_index"
$ws.Range("B64").Value = 0
$ws.Range("A65").Value = "param_P_pv2_bat2_starting_index"
$ws.Range("B65").Value = 0
$ws.Range("A66").Value = "param_P_pv2_demand2_starting_index"
$ws.Range("B66").Value = 0
$ws.Range("A67").Value = "param_P_pv2_heat_pump1_starting_index"
$ws.Range("B67").Value = 0
$ws.Range("A68").Value = "param_P_pv2_heat_pump2_starting_index"
$ws.Range("B68").Value = 0
$ws.Range("A69").Value = "param_pv2_inv_cost_starting_index"
$ws.Range("B69").Value = 0
$ws.Range("A70").Value = "param_P_pv2_charging_station2_starting_index"
$ws.Range("B70").Value = 0
$ws.Range("A71").Value = "param_P_pv2_demand1_starting_index"
$ws.Range("B71").Value = 1.8
$ws.Range("A72").Value = "param_P_pv2_net2_starting_index"
$ws.Range("B72").Value = 0
$ws.Range("A73").Value = "param_pv2_op_cost_starting_index"
$ws.Range("B73").Value = 1
$ws.Range("A74").Value = "param_pv2_emissions_starting_index"
$ws.Range("B74").Value = 0.8999999999999999
$ws.Range("A75").Value = "param_bat1_K_ch_starting_index"
$ws.Range("B75").Value = 0
$ws.Range("A76").Value = "param_P_bat1_net1_starting_index"
$ws.Range("B76").Value = 0
$ws.Range("A77").Value = "param_P_bat1_heat_pump1_starting_index"
$ws.Range("B77").Value = 0
$ws.Range("A78").Value = "param_P_bat1_demand2_starting_index"
$ws.Range("B78").Value = 0
$ws.Range("A79").Value = "param_bat1_cumulated_aging_starting_index"
$ws.Range("B79").Value = 0.000003166666666666667
$ws.Range("A80").Value = "param_bat1_op_cost_starting_index"
$ws.Range("B80").Value = 1
$ws.Range("A81").Value = "param_bat1_emissions_starting_index"
$ws.Range("B81").Value = 0
$ws.Range("A82").Value = "param_bat1_K_dis_starting_index"
$ws.Range("B82").Value = 1
$ws.Range("A83").Value = "param_bat1_integer_starting_index"
$ws.Range("B83").Value = -0
$ws.Range("A84").Value = "param_bat1_SOC_starting_index"
$ws.Range("B84").Value = 0.3
$ws.Range("A85").Value = "param_P_bat1_charging_station1_starting_index"
$ws.Range("B85").Value = 0
$ws.Range("A86").Value = "param_P_bat1_heat_pump2_starting_index"
$ws.Range("B86").Value = 0
$ws.Range("A87").Value = "param_P_to_bat1_starting_index"
$ws.Range("B87").Value = 0
$ws.Range("A88").Value = "param_bat1_inv_cost_starting_index"
$ws.Range("B88").Value = 0
$ws.Range("A89").Value = "param_P_from_bat1_starting_index"
$ws.Range("B89").Value = 0
$ws.Range("A90").Value = "param_P_bat1_demand1_starting_index"
$ws.Range("B90").Value = 0
$ws.Range("A91").Value = "param_P_bat1_charging_station2_starting_index"
$ws.Range("B91").Value = 0
$ws.Range("A92").Value = "param_P_bat1_net2_starting_index"
$ws.Range("B92").Value = 0
$ws.Range("A93").Value = "param_bat1_SOC_max_starting_index"
$ws.Range("B93").Value = 0.9999968333333333
$ws.Range("A94").Value = "param_P_bat2_net2_starting_index"
$ws.Range("B94").Value = 0
$ws.Range("A95").Value = "param_P_bat2_demand1_starting_index"
$ws.Range("B95").Value = 0
$ws.Range("A96").Value = "param_bat2_SOC_starting_index"
$ws.Range("B96").Value = 0.5
$ws.Range("A97").Value = "param_P_bat2_charging_station1_starting_index"
$ws.Range("B97").Value = 0
$ws.Range("A98").Value = "param_bat2_op_cost_starting_index"
$ws.Range("B98").Value = 1
$ws.Range("A99").Value = "param_P_bat2_net1_starting_index"
$ws.Range("B99").Value = 0
$ws.Range("A100").Value = "param_bat2_inv_cost_starting_index"
$ws.Range("B100").Value = 0
$ws.Range("A101").Value = "param_bat2_emissions_starting_index"
$ws.Range("B101").Value = 0
$ws.Range("A102").Value = "param_bat2_cumulated_aging_starting_index"
$ws.Range("B102").Value = 0
$ws.Range("A103").Value = "param_P_bat2_demand2_starting_index"
$ws.Range("B103").Value = -0
$ws.Range("A104").Value = "param_P_to_bat2_starting_index"
$ws.Range("B104").Value = 0
$ws.Range("A105").Value = "param_P_bat2_charging_station2_starting_index"
$ws.Range("B105").Value = 0
$ws.Range("A106").Value = "param_P_bat2_heat_pump1_starting_index"
$ws.Range("B106").Value = 0
$ws.Range("A107").Value = "param_P_bat2_heat_pump2_starting_index"
$ws.Range("B107").Value = 0
$ws.Range("A108").Value = "param_P_from_bat2_starting_index"
$ws.Range("B108").Value = 0
$ws.Range("A109").Value = "param_bat2_K_dis_starting_index"
$ws.Range("B109").Value = 1
$ws.Range("A110").Value = "param_bat2_K_ch_starting_index"
$ws.Range("B110").Value = 0
$ws.Range("A111").Value = "param_bat2_SOC_max_starting_index"
$ws.Range("B111").Value = 1
$ws.Range("A112").Value = "param_bat2_integer_starting_index"
$ws.Range("B112").Value = -0
$ws.Range("A113").Value = "param_Q_CHP1_demand1_starting_index"
$ws.Range("B113").Value = 40
$ws.Range("A114").Value = "param_P_CHP1_bat2_starting_index"
$ws.Range("B114").Value = 0
$ws.Range("A115").Value = "param_P_CHP1_charging_station2_starting_index"
$ws.Range("B115").Value = 0
$ws.Range("A116").Value = "param_P_CHP1_demand2_starting_index"
$ws.Range("B116").Value = 0
$ws.Range("A117").Value = "param_CHP1_emissions_starting_index"
$ws.Range("B117").Value = 4.83
$ws.Range("A118").Value = "param_P_CHP1_heat_pump1_starting_index"
$ws.Range("B118").Value = 0
$ws.Range("A119").Value = "param_P_CHP1_bat1_starting_index"
$ws.Range("B119").Value = 0
$ws.Range("A120").Value = "param_P_CHP1_net2_starting_index"
$ws.Range("B120").Value = 0
$ws.Range("A121").Value = "param_CHP1_inv_cost_starting_index"
$ws.Range("B121").Value = 0
$ws.Range("A122").Value = "param_P_CHP1_net1_starting_index"
$ws.Range("B122").Value = 0
$ws.Range("A123").Value = "param_Q_CHP1_net2_starting_index"
$ws.Range("B123").Value = 0
$ws.Range("A124").Value = "param_P_CHP1_charging_station1_starting_index"
$ws.Range("B124").Value = 20
$ws.Range("A125").Value = "param_P_CHP1_heat_pump2_starting_index"
$ws.Range("B125").Value = 0
$ws.Range("A126").Value = "param_P_from_CHP1_starting_index"
$ws.Range("B126").Value = 20
$ws.Range("A127").Value = "param_Q_CHP1_net1_starting_index"
$ws.Range("B127").Value = 0
$ws.Range("A128").Value = "param_P_CHP1_demand1_starting_index"
$ws.Range("B128").Value = 0
$ws.Range("A129").Value = "param_CHP1_fuel_cons_starting_index"
$ws.Range("B129").Value = 2.1
$ws.Range("A130").Value = "param_CHP1_op_cost_starting_index"
$ws.Range("B130").Value = 10.5
$ws.Range("A131").Value = "param_Q_from_CHP1_starting_index"
$ws.Range("B131").Value = 40
$ws.Range("A132").Value = "param_Q_CHP1_demand2_starting_index"
$ws.Range("B132").Value = 0
$ws.Range("A133").Value = "param_P_CHP2_net2_starting_index"
$ws.Range("B133").Value = 0
$ws.Range("A134").Value = "param_P_CHP2_bat2_starting_index"
$ws.Range("B134").Value = 0
$ws.Range("A135").Value = "param_P_from_CHP2_starting_index"
$ws.Range("B135").Value = 20
$ws.Range("A136").Value = "param_CHP2_op_cost_starting_index"
$ws.Range("B136").Value = 10.5
$ws.Range("A137").Value = "param_Q_CHP2_net1_starting_index"
$ws.Range("B137").Value = 0
$ws.Range("A138").Value = "param_P_CHP2_heat_pump2_starting_index"
$ws.Range("B138").Value = 0
$ws.Range("A139").Value = "param_Q_CHP2_net2_starting_index"
$ws.Range("B139").Value = 0
$ws.Range("A140").Value = "param_Q_CHP2_demand1_starting_index"
$ws.Range("B140").Value = 40
$ws.Range("A141").Value = "param_CHP2_fuel_cons_starting_index"
$ws.Range("B141").Value = 2.1
$ws.Range("A142").Value = "param_P_CHP2_charging_station1_starting_index"
$ws.Range("B142").Value = 20
$ws.Range("A143").Value = "param_CHP2_inv_cost_starting_index"
$ws.Range("B143").Value = 0
$ws.Range("A144").Value = "param_P_CHP2_bat1_starting_index"
$ws.Range("B144").Value = 0
$ws.Range("A145").Value = "param_P_CHP2_demand1_starting_index"
$ws.Range("B145").Value = 0
$ws.Range("A146").Value = "param_Q_from_CHP2_starting_index"
$ws.Range("B146").Value = 40
$ws.Range("A147").Value = "param_Q_CHP2_demand2_starting_index"
$ws.Range("B147").Value = 0
$ws.Range("A148").Value = "param_P_CHP2_heat_pump1_starting_index"
$ws.Range("B148").Value = 0
$ws.Range("A149").Value = "param_P_CHP2_net1_starting_index"
$ws.Range("B149").Value = 0
$ws.Range("A150").Value = "param_CHP2_emissions_starting_index"
$ws.Range("B150").Value = 4.83
$ws.Range("A151").Value = "param_P_CHP2_demand2_starting_index"
$ws.Range("B151").Value = 0
$ws.Range("A152").Value = "param_P_CHP2_charging_station2_starting_index"
$ws.Range("B152").Value = 0
$ws.Range("A153").Value = "param_Q_solar_th1_net2_starting_index"
$ws.Range("B153").Value = 0
$ws.Range("A154").Value = "param_Q_from_solar_th1_starting_index"
$ws.Range("B154").Value = 0.1111112
$ws.Range("A155").Value = "param_Q_solar_th1_demand1_starting_index"
$ws.Range("B155").Value = 0.1111112
$ws.Range("A156").Value = "param_Q_solar_th1_net1_starting_index"
$ws.Range("B156").Value = 0
$ws.Range("A157").Value = "param_solar_th1_op_cost_starting_index"
$ws.Range("B157").Value = 1
$ws.Range("A158").Value = "param_solar_th1_inv_cost_starting_index"
$ws.Range("B158").Value = 0
$ws.Range("A159").Value = "param_solar_th1_emissions_starting_index"
$ws.Range("B159").Value = 0.05555560000000001
$ws.Range("A160").Value = "param_Q_solar_th1_demand2_starting_index"
$ws.Range("B160").Value = 0
$ws.Range("A161").Value = "param_Q_solar_th2_net1_starting_index"
$ws.Range("B161").Value = 0
$ws.Range("A162").Value = "param_Q_solar_th2_net2_starting_index"
$ws.Range("B162").Value = 0
$ws.Range("A163").Value = "param_Q_solar_th2_demand2_starting_index"
$ws.Range("B163").Value = 0
$ws.Range("A164").Value = "param_solar_th2_emissions_starting_index"
$ws.Range("B164").Value = 0.6000000000000001
$ws.Range("A165").Value = "param_Q_from_solar_th2_starting_index"
$ws.Range("B165").Value = 1.2
$ws.Range("A166").Value = "param_solar_th2_op_cost_starting_index"
$ws.Range("B166").Value = 1
$ws.Range("A167").Value = "param_Q_solar_th2_demand1_starting_index"
$ws.Range("B167").Value = 1.2
$ws.Range("A168").Value = "param_solar_th2_inv_cost_starting_index"
$ws.Range("B168").Value = 0
$ws.Range("A169").Value = "param_P_pvt1_bat2_starting_index"
$ws.Range("B169").Value = 0
$ws.Range("A170").Value = "param_P_from_pvt1_starting_index"
$ws.Range("B170").Value = 0.2222224
$ws.Range("A171").Value = "param_pvt1_emissions_starting_index"
$ws.Range("B171").Value = 0.14444456
$ws.Range("A172").Value = "param_pvt1_inv_cost_starting_index"
$ws.Range("B172").Value = 0
$ws.Range("A173").Value = "param_P_pvt1_bat1_starting_index"
$ws.Range("B173").Value = 0
$ws.Range("A174").Value = "param_P_pvt1_net2_starting_index"
$ws.Range("B174").Value = 0
$ws.Range("A175").Value = "param_pvt1_op_cost_starting_index"
$ws.Range("B175").Value = 1
$ws.Range("A176").Value = "param_P_pvt1_heat_pump2_starting_index"
$ws.Range("B176").Value = 0.1111112
$ws.Range("A177").Value = "param_P_pvt1_charging_station2_starting_index"
$ws.Range("B177").Value = 0
$ws.Range("A178").Value = "param_Q_pvt1_net1_starting_index"
$ws.Range("B178").Value = 0
$ws.Range("A179").Value = "param_Q_from_pvt1_starting_index"
$ws.Range("B179").Value = 0.2888891200000001
$ws.Range("A180").Value = "param_P_pvt1_net1_starting_index"
$ws.Range("B180").Value = 0
$ws.Range("A181").Value = "param_Q_pvt1_demand2_starting_index"
$ws.Range("B181").Value = 0
$ws.Range("A182").Value = "param_P_pvt1_demand1_starting_index"
$ws.Range("B182").Value = 0
$ws.Range("A183").Value = "param_P_pvt1_heat_pump1_starting_index"
$ws.Range("B183").Value = 0.1111112
$ws.Range("A184").Value = "param_Q_pvt1_net2_starting_index"
$ws.Range("B184").Value = 0
$ws.Range("A185").Value = "param_P_pvt1_charging_station1_starting_index"
$ws.Range("B185").Value = 0
$ws.Range("A186").Value = "param_Q_pvt1_demand1_starting_index"
$ws.Range("B186").Value = 0.2888891200000001
$ws.Range("A187").Value = "param_P_pvt1_demand2_starting_index"
$ws.Range("B187").Value = 0
$ws.Range("A188").Value = "param_P_pvt2_bat1_starting_index"
$ws.Range("B188").Value = 0
$ws.Range("A189").Value = "param_P_pvt2_demand2_starting_index"
$ws.Range("B189").Value = 0
$ws.Range("A190").Value = "param_P_pvt2_net1_starting_index"
$ws.Range("B190").Value = 0
$ws.Range("A191").Value = "param_Q_from_pvt2_starting_index"
$ws.Range("B191").Value = 1.56
$ws.Range("A192").Value = "param_P_from_pvt2_starting_index"
$ws.Range("B192").Value = 1.2
$ws.Range("A193").Value = "param_P_pvt2_charging_station2_starting_index"
$ws.Range("B193").Value = 0
$ws.Range("A194").Value = "param_P_pvt2_heat_pump2_starting_index"
$ws.Range("B194").Value = 0
$ws.Range("A195").Value = "param_P_pvt2_bat2_starting_index"
$ws.Range("B195").Value = 0
$ws.Range("A196").Value = "param_pvt2_emissions_starting_index"
$ws.Range("B196").Value = 0.7800000000000001
$ws.Range("A197").Value = "param_Q_pvt2_demand1_starting_index"
$ws.Range("B197").Value = 1.56
$ws.Range("A198").Value = "param_P_pvt2_charging_station1_starting_index"
$ws.Range("B198").Value = 0
$ws.Range("A199").Value = "param_P_pvt2_demand1_starting_index"
$ws.Range("B199").Value = 1.2
$ws.Range("A200").Value = "param_Q_pvt2_net2_starting_index"
$ws.Range("B200").Value = 0
$ws.Range("A201").Value = "param_P_pvt2_net2_starting_index"
$ws.Range("B201").Value = 0
$ws.Range("A202").Value = "param_P_pvt2_heat_pump1_starting_index"
$ws.Range("B202").Value = 0
$ws.Range("A203").Value = "param_Q_pvt2_net1_starting_index"
$ws.Range("B203").Value = 0
$ws.Range("A204").Value = "param_pvt2_inv_cost_starting_index"
$ws.Range("B204").Value = 0
$ws.Range("A205").Value = "param_Q_pvt2_demand2_starting_index"
$ws.Range("B205").Value = 0
$ws.Range("A206").Value = "param_pvt2_op_cost_starting_index"
$ws.Range("B206").Value = 1
$ws.Range("A207").Value = "param_charging_station1_inv_cost_starting_index"
$ws.Range("B207").Value = 0
$ws.Range("A208").Value = "param_charging_station1_op_cost_starting_index"
$ws.Range("B208").Value = -29.96042500000001
$ws.Range("A209").Value = "param_charging_station1_emissions_starting_index"
$ws.Range("B209").Value = 2.496702083333334
$ws.Range("A210").Value = "param_charging_station2_emissions_starting_index"
$ws.Range("B210").Value = 1.65
$ws.Range("A211").Value = "param_charging_station2_inv_cost_starting_index"
$ws.Range("B211").Value = 0
$ws.Range("A212").Value = "param_charging_station2_op_cost_starting_index"
$ws.Range("B212").Value = -19.8
$ws.Range("A213").Value = "param_Q_heat_pump1_net1_starting_index"
$ws.Range("B213").Value = 0
$ws.Range("A214").Value = "param_heat_pump1_op_cost_starting_index"
$ws.Range("B214").Value = 8.561643835616438
$ws.Range("A215").Value = "param_heat_pump1_emissions_starting_index"
$ws.Range("B215").Value = 2.76
$ws.Range("A216").Value = "param_Q_from_heat_pump1_starting_index"
$ws.Range("B216").Value = 80
$ws.Range("A217").Value = "param_Q_to_heat_pump1_starting_index"
$ws.Range("B217").Value = 0
$ws.Range("A218").Value = "param_P_from_heat_pump1_starting_index"
$ws.Range("B218").Value = 0
$ws.Range("A219").Value = "param_heat_pump1_inv_cost_starting_index"
$ws.Range("B219").Value = 0
$ws.Range("A220").Value = "param_Q_heat_pump1_net2_starting_index"
$ws.Range("B220").Value = 0
$ws.Range("A221").Value = "param_Q_heat_pump1_demand1_starting_index"
$ws.Range("B221").Value = 80
$ws.Range("A222").Value = "param_Q_heat_pump1_demand2_starting_index"
$ws.Range("B222").Value = 0
$ws.Range("A223").Value = "param_P_to_heat_pump1_starting_index"
$ws.Range("B223").Value = 20
$ws.Range("A224").Value = "param_P_to_heat_pump2_starting_index"
$ws.Range("B224").Value = 20
$ws.Range("A225").Value = "param_Q_from_heat_pump2_starting_index"
$ws.Range("B225").Value = 80
$ws.Range("A226").Value = "param_heat_pump2_emissions_starting_index"
$ws.Range("B226").Value = 2.76
$ws.Range("A227").Value = "param_P_from_heat_pump2_starting_index"
$ws.Range("B227").Value = 0
$ws.Range("A228").Value = "param_Q_heat_pump2_demand2_starting_index"
$ws.Range("B228").Value = 0
$ws.Range("A229").Value = "param_Q_heat_pump2_net1_starting_index"
$ws.Range("B229").Value = 0
$ws.Range("A230").Value = "param_Q_heat_pump2_net2_starting_index"
$ws.Range("B230").Value = 0
$ws.Range("A231").Value = "param_Q_heat_pump2_demand1_starting_index"
$ws.Range("B231").Value = 80
$ws.Range("A232").Value = "param_heat_pump2_op_cost_starting_index"
$ws.Range("B232").Value = 8.561643835616438
$ws.Range("A233").Value = "param_heat_pump2_inv_cost_starting_index"
$ws.Range("B233").Value = 0
$ws.Range("A234").Value = "param_Q_to_heat_pump2_starting_index"
$ws.Range("B234").Value = 0
$ws.Range("A235").Value = "param_total_emissions_starting_index"
$ws.Range("B235").Value = 852.3124636167672
$ws.Range("A236").Value = "param_total_sell_starting_index"
$ws.Range("B236").Value = 0
$ws.Range("A237").Value = "param_total_buy_starting_index"
$ws.Range("B237").Value = 668.9583166453599
$ws.Range("A238").Value = "param_total_operation_cost_starting_index"
$ws.Range("B238").Value = -7.637137328767139
